# Re-balance the jh/sh candidate groups from 5 groups of ~4 down to
# 4 larger groups (J01-J04 / S01-S04 instead of J01-J05 / S01-S05),
# folding the last (smaller) group into the others.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("jh_candidate_data")
$ws1.Range("C6").Value  = "J01"
$ws1.Range("C10").Value = "J02"
$ws1.Range("C11").Value = "J02"
$ws1.Range("C14").Value = "J03"
$ws1.Range("C15").Value = "J03"
$ws1.Range("C17").Value = "J04"
$ws1.Range("C18").Value = "J04"
$ws1.Range("C19").Value = "J04"

$ws2 = $wb.Worksheets.Item("sh_candidate_data")
$ws2.Range("C6").Value  = "S01"
$ws2.Range("C10").Value = "S02"
$ws2.Range("C11").Value = "S02"
$ws2.Range("C14").Value = "S03"
$ws2.Range("C15").Value = "S03"
$ws2.Range("C17").Value = "S04"
$ws2.Range("C18").Value = "S04"
$ws2.Range("C19").Value = "S04"

# Leave the selection/active-sheet state matching where the author's
# cursor ended up: jh sheet group_id column selected first, then the
# sh sheet's last group selected last (making sh_candidate_data the
# active tab on save).
$ws1.Range("C2:C19").Select()
$ws2.Range("C17:C19").Select()
